# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Macroferia Regional de Talca - Piña"
# at row 409, shifting the existing rows 409..490 down to 410..491.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 409 (pushes rows 409-490 down to 410-491)
$ws.Rows.Item(409).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Cells.Item(409, 1).Value  = 5
$ws.Cells.Item(409, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(409, 3).Value  = "Maule"
$ws.Cells.Item(409, 4).Value  = 45258
$ws.Cells.Item(409, 5).Value  = 7
$ws.Cells.Item(409, 6).Value  = "Fruta"
$ws.Cells.Item(409, 7).Value  = 100108
$ws.Cells.Item(409, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(409, 9).Value  = 100108005
$ws.Cells.Item(409, 10).Value = "Piña"
$ws.Cells.Item(409, 11).Value = "Caramelo"
$ws.Cells.Item(409, 12).Value = "Segunda"
$ws.Cells.Item(409, 13).Value = 120
$ws.Cells.Item(409, 14).Value = 22000
$ws.Cells.Item(409, 15).Value = 22000
$ws.Cells.Item(409, 16).Value = 22000
$ws.Cells.Item(409, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(409, 18).Value = "Ecuador"
$ws.Cells.Item(409, 19).Value = 1571
$ws.Cells.Item(409, 20).Value = 14
